# Weekly update for "Hortaliza, Macroferia Regional de Talca - Espárragos".
# Two new weekly records are inserted into the price table:
#   - one at row 51 (pushing the former rows 51-53 down to 52-54)
#   - another at row 55 (pushing the former rows 54-70, now at 54-56..72,
#     further down by one more)
# so the sheet grows from A1:R70 to A1:R72.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- First insertion: new row at position 51 ---------------------------
$ws.Rows(51).Insert()

$ws.Cells.Item(51, 1).Value = 5
$ws.Cells.Item(51, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(51, 3).Value = "Maule"
$ws.Cells.Item(51, 4).Value = 44846
$ws.Cells.Item(51, 4).NumberFormat = $ws.Cells.Item(52, 4).NumberFormat
$ws.Cells.Item(51, 5).Value = 7
$ws.Cells.Item(51, 6).Value = 300000000
$ws.Cells.Item(51, 7).Value = "Espárragos"
$ws.Cells.Item(51, 8).Value = "Sin especificar"
$ws.Cells.Item(51, 9).Value = "Primera"
$ws.Cells.Item(51, 10).Value = 3000
$ws.Cells.Item(51, 11).Value = 1200
$ws.Cells.Item(51, 12).Value = 1200
$ws.Cells.Item(51, 13).Value = 1200
$ws.Cells.Item(51, 14).Value = "$/kilo"
$ws.Cells.Item(51, 15).Value = "Provincia de Linares"
$ws.Cells.Item(51, 16).Value = 1200
$ws.Cells.Item(51, 17).Value = 1
$ws.Cells.Item(51, 18).Value = "Hortaliza"

# --- Second insertion: new row at position 55 ---------------------------
$ws.Rows(55).Insert()

$ws.Cells.Item(55, 1).Value = 5
$ws.Cells.Item(55, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(55, 3).Value = "Maule"
$ws.Cells.Item(55, 4).Value = 44845
$ws.Cells.Item(55, 4).NumberFormat = $ws.Cells.Item(54, 4).NumberFormat
$ws.Cells.Item(55, 5).Value = 7
$ws.Cells.Item(55, 6).Value = 300000000
$ws.Cells.Item(55, 7).Value = "Espárragos"
$ws.Cells.Item(55, 8).Value = "Sin especificar"
$ws.Cells.Item(55, 9).Value = "Primera"
$ws.Cells.Item(55, 10).Value = 3000
$ws.Cells.Item(55, 11).Value = 1200
$ws.Cells.Item(55, 12).Value = 1200
$ws.Cells.Item(55, 13).Value = 1200
$ws.Cells.Item(55, 14).Value = "$/kilo"
$ws.Cells.Item(55, 15).Value = "Provincia de Linares"
$ws.Cells.Item(55, 16).Value = 1200
$ws.Cells.Item(55, 17).Value = 1
$ws.Cells.Item(55, 18).Value = "Hortaliza"
